$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.95625533333333
$ws.Range("H2").Value = 53.86876599999999
$ws.Range("I2").Value = 0.05960074617816258
$ws.Range("J2").Value = 0.05960074617816258
$ws.Range("M2").Value = 13.17295566666667
$ws.Range("N2").Value = 39.518867
$ws.Range("O2").Value = 0.133784132206724
$ws.Range("P2").Value = 0.133784132206724
$ws.Range("Q2").Value = 236.5369554453469
$ws.Range("R2").Value = 2128.832599008122
$ws.Range("S2").Value = 0.007973634106318702
$ws.Range("T2").Value = 0.007973634106318704
$ws.Range("G3").Value = 17.95625533333333
$ws.Range("H3").Value = 53.86876599999999
$ws.Range("I3").Value = 0.05960074617816258
$ws.Range("J3").Value = 0.05960074617816258
$ws.Range("O3").Value = 0.4382627974978752
$ws.Range("P3").Value = 0.4382627974978752
$ws.Range("Q3").Value = 774.8702786734357
$ws.Range("R3").Value = 6973.832508060921
$ws.Range("S3").Value = 0.02612078975300233
$ws.Range("T3").Value = 0.02612078975300233
$ws.Range("G4").Value = 17.95625533333333
$ws.Range("H4").Value = 53.86876599999999
$ws.Range("I4").Value = 0.05960074617816258
$ws.Range("J4").Value = 0.05960074617816258
$ws.Range("M4").Value = 21.06166566666667
$ws.Range("N4").Value = 63.184997
$ws.Range("O4").Value = 0.2139016281041017
$ws.Range("P4").Value = 0.2139016281041017
$ws.Range("Q4").Value = 378.1886464559668
$ws.Range("R4").Value = 3403.697818103702
$ws.Range("S4").Value = 0.0127486966437283
$ws.Range("T4").Value = 0.0127486966437283
$ws.Range("G5").Value = 17.95625533333333
$ws.Range("H5").Value = 53.86876599999999
$ws.Range("I5").Value = 0.05960074617816258
$ws.Range("J5").Value = 0.05960074617816258
$ws.Range("M5").Value = 21.076417
$ws.Range("N5").Value = 63.229251
$ws.Range("O5").Value = 0.214051442191299
$ws.Range("P5").Value = 0.214051442191299
$ws.Range("Q5").Value = 378.4535251638073
$ws.Range("R5").Value = 3406.081726474266
$ws.Range("S5").Value = 0.01275762567511326
$ws.Range("T5").Value = 0.01275762567511326
$ws.Range("H6").Value = 789.271408
$ws.Range("I6").Value = 0.8732549183303921
$ws.Range("J6").Value = 0.8732549183303921
$ws.Range("M6").Value = 13.17295566666667
$ws.Range("N6").Value = 39.518867
$ws.Range("O6").Value = 0.133784132206724
$ws.Range("P6").Value = 0.133784132206724
$ws.Range("Q6").Value = 3465.679088850526
$ws.Range("R6").Value = 31191.11179965473
$ws.Range("S6").Value = 0.1168276514440851
$ws.Range("T6").Value = 0.1168276514440851
$ws.Range("H7").Value = 789.271408
$ws.Range("I7").Value = 0.8732549183303921
$ws.Range("J7").Value = 0.8732549183303921
$ws.Range("O7").Value = 0.4382627974978752
$ws.Range("P7").Value = 0.4382627974978752
$ws.Range("S7").Value = 0.3827151434362562
$ws.Range("T7").Value = 0.3827151434362562
$ws.Range("H8").Value = 789.271408
$ws.Range("I8").Value = 0.8732549183303921
$ws.Range("J8").Value = 0.8732549183303921
$ws.Range("M8").Value = 21.06166566666667
$ws.Range("N8").Value = 63.184997
$ws.Range("O8").Value = 0.2139016281041017
$ws.Range("P8").Value = 0.2139016281041017
$ws.Range("Q8").Value = 5541.123505185085
$ws.Range("R8").Value = 49870.11154666577
$ws.Range("S8").Value = 0.1867906487807853
$ws.Range("T8").Value = 0.1867906487807853
$ws.Range("H9").Value = 789.271408
$ws.Range("I9").Value = 0.8732549183303921
$ws.Range("J9").Value = 0.8732549183303921
$ws.Range("M9").Value = 21.076417
$ws.Range("N9").Value = 63.229251
$ws.Range("O9").Value = 0.214051442191299
$ws.Range("P9").Value = 0.214051442191299
$ws.Range("Q9").Value = 5545.004440395046
$ws.Range("R9").Value = 49905.03996355541
$ws.Range("S9").Value = 0.1869214746692655
$ws.Range("T9").Value = 0.1869214746692655
$ws.Range("G10").Value = 13.13303333333333
$ws.Range("H10").Value = 39.3991
$ws.Range("I10").Value = 0.04359141545488614
$ws.Range("J10").Value = 0.04359141545488615
$ws.Range("M10").Value = 13.17295566666667
$ws.Range("N10").Value = 39.518867
$ws.Range("O10").Value = 0.133784132206724
$ws.Range("P10").Value = 0.133784132206724
$ws.Range("Q10").Value = 173.0008658688555
$ws.Range("R10").Value = 1557.0077928197
$ws.Range("S10").Value = 0.005831839688294718
$ws.Range("T10").Value = 0.005831839688294721
$ws.Range("G11").Value = 13.13303333333333
$ws.Range("H11").Value = 39.3991
$ws.Range("I11").Value = 0.04359141545488614
$ws.Range("J11").Value = 0.04359141545488615
$ws.Range("O11").Value = 0.4382627974978752
$ws.Range("P11").Value = 0.4382627974978752
$ws.Range("Q11").Value = 566.732707344411
$ws.Range("R11").Value = 5100.594366099699
$ws.Range("S11").Value = 0.01910449568415051
$ws.Range("T11").Value = 0.01910449568415052
$ws.Range("G12").Value = 13.13303333333333
$ws.Range("H12").Value = 39.3991
$ws.Range("I12").Value = 0.04359141545488614
$ws.Range("J12").Value = 0.04359141545488615
$ws.Range("M12").Value = 21.06166566666667
$ws.Range("N12").Value = 63.184997
$ws.Range("O12").Value = 0.2139016281041017
$ws.Range("P12").Value = 0.2139016281041017
$ws.Range("Q12").Value = 276.6035572558555
$ws.Range("R12").Value = 2489.4320153027
$ws.Range("S12").Value = 0.00932427473716245
$ws.Range("T12").Value = 0.00932427473716245
$ws.Range("G13").Value = 13.13303333333333
$ws.Range("H13").Value = 39.3991
$ws.Range("I13").Value = 0.04359141545488614
$ws.Range("J13").Value = 0.04359141545488615
$ws.Range("M13").Value = 21.076417
$ws.Range("N13").Value = 63.229251
$ws.Range("O13").Value = 0.214051442191299
$ws.Range("P13").Value = 0.214051442191299
$ws.Range("Q13").Value = 276.7972870082334
$ws.Range("R13").Value = 2491.1755830741
$ws.Range("S13").Value = 0.009330805345278461
$ws.Range("T13").Value = 0.009330805345278463
$ws.Range("G14").Value = 7.095921999999999
$ws.Range("H14").Value = 21.287766
$ws.Range("I14").Value = 0.02355292003655921
$ws.Range("J14").Value = 0.02355292003655921
$ws.Range("M14").Value = 13.17295566666667
$ws.Range("N14").Value = 39.518867
$ws.Range("O14").Value = 0.133784132206724
$ws.Range("P14").Value = 0.133784132206724
$ws.Range("Q14").Value = 93.47426592012465
$ws.Range("R14").Value = 841.2683932811219
$ws.Range("S14").Value = 0.003151006968025435
$ws.Range("T14").Value = 0.003151006968025436
$ws.Range("G15").Value = 7.095921999999999
$ws.Range("H15").Value = 21.287766
$ws.Range("I15").Value = 0.02355292003655921
$ws.Range("J15").Value = 0.02355292003655921
$ws.Range("O15").Value = 0.4382627974978752
$ws.Range("P15").Value = 0.4382627974978752
$ws.Range("Q15").Value = 306.2118997259913
$ws.Range("R15").Value = 2755.907097533922
$ws.Range("S15").Value = 0.01032236862446619
$ws.Range("T15").Value = 0.01032236862446619
$ws.Range("G16").Value = 7.095921999999999
$ws.Range("H16").Value = 21.287766
$ws.Range("I16").Value = 0.02355292003655921
$ws.Range("J16").Value = 0.02355292003655921
$ws.Range("M16").Value = 21.06166566666667
$ws.Range("N16").Value = 63.184997
$ws.Range("O16").Value = 0.2139016281041017
$ws.Range("P16").Value = 0.2139016281041017
$ws.Range("Q16").Value = 149.4519367607446
$ws.Range("R16").Value = 1345.067430846702
$ws.Range("S16").Value = 0.005038007942425734
$ws.Range("T16").Value = 0.005038007942425734
$ws.Range("G17").Value = 7.095921999999999
$ws.Range("H17").Value = 21.287766
$ws.Range("I17").Value = 0.02355292003655921
$ws.Range("J17").Value = 0.02355292003655921
$ws.Range("M17").Value = 21.076417
$ws.Range("N17").Value = 63.229251
$ws.Range("O17").Value = 0.214051442191299
$ws.Range("P17").Value = 0.214051442191299
$ws.Range("Q17").Value = 149.556611071474
$ws.Range("R17").Value = 1346.009499643266
$ws.Range("S17").Value = 0.005041536501641842
$ws.Range("T17").Value = 0.005041536501641842
